# Apply the changes described by the diff:
#  - shared string "Timeseries" -> "Annual" (used in cell B7, merged B7:B10)
#  - sheet selection changed from C12 to B7:B10

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("temp_metric_table")

# Update the cell text that currently holds "Timeseries"
$ws.Range("B7").Value = "Annual"

# Update the active selection to B7:B10
$ws.Range("B7:B10").Select()
